# Update "Hjemme passive" meanEMG / legmaxROM data: subjects 1-4 (cols B-E)
# replaced with data for subjects 15,16,15,16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - subject id header
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - CON
$ws.Range("B2").Value = 285.29093073464242
$ws.Range("C2").Value = 265.6287409710601
$ws.Range("D2").Value = 284.5956040154266
$ws.Range("E2").Value = 262.21989162037539

# Row 3 - STR
$ws.Range("B3").Value = 308.9846248534343
$ws.Range("C3").Value = 261.16854873030132
$ws.Range("D3").Value = 318.09569398032079
$ws.Range("E3").Value = 261.79447028965973

# Update the selection to match the new narrower region of interest
$ws.Range("B1:E3").Select()
